$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# H1: header label changed from "fess" to "0010" (kept as text, not numeric)
$ws.Range("H1").Formula = "'0010"

# H2 and H3: formula-error (#NAME?) results replaced with plain numeric values
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 2

# H4:H31: formula-error (#NAME?) results cleared out to blank cells
$ws.Range("H4:H31").ClearContents()
